# Atualização automática via cronjob
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (values that changed in-place) ---
$ws.Range("A2").Value  = 8
$ws.Range("A3").Value  = 12
$ws.Range("A7").Value  = 7
$ws.Range("A8").Value  = 9
$ws.Range("A9").Value  = 10
$ws.Range("A10").Value = 13
$ws.Range("A11").Value = 11

$ws.Range("C12").Value = 400

# Row 13 is repurposed with new data (old row 13 content moves to the new row 17 below)
$ws.Range("A13").Value = 4
$ws.Range("C13").Value = 400
$ws.Range("E13").Value = "000787"
$ws.Range("F13").Value = "SACO DE LIXO 50L COMUM PACOTINHO C/10 UND"
$ws.Range("G13").Value = -88

# --- Append new rows 14-17 ---
# First clone formatting (incl. the bordered/centered style used in column A)
# from row 13 down onto rows 14-17, then set the actual values.
$ws.Range("A13:H13").Copy($ws.Range("A14"))
$ws.Range("A13:H13").Copy($ws.Range("A15"))
$ws.Range("A13:H13").Copy($ws.Range("A16"))
$ws.Range("A13:H13").Copy($ws.Range("A17"))

# Row 14
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "2025-04-10"
$ws.Range("C14").Value = 400
$ws.Range("D14").Value = "MAYCA CONSTRUCOES E INSTALACOES ELETROMECANICAS EIRELI"
$ws.Range("E14").Value = "000188"
$ws.Range("F14").Value = "SACO DE LIXO 100L COMUM PACOTINHO C/5 UND"
$ws.Range("G14").Value = 623
$ws.Range("H14").Value = $false

# Row 15
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "2025-04-10"
$ws.Range("C15").Value = 28
$ws.Range("D15").Value = "LABORATORIO SANTOS E VIDAL LTDA"
$ws.Range("E15").Value = "000897"
$ws.Range("F15").Value = "AGUA SANITARIA GLOBO SAN 5L"
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = $false

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "2025-04-10"
$ws.Range("C16").Value = 80
$ws.Range("D16").Value = "DE PASQUAL HOTEIS E TURISMO LTDA."
$ws.Range("E16").Value = "010189"
$ws.Range("F16").Value = "PILHA ALCALINA AAA PANASONIC"
$ws.Range("G16").Value = 112
$ws.Range("H16").Value = $false

# Row 17 (this is the content that used to live in row 13)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "2025-04-10"
$ws.Range("C17").Value = 200
$ws.Range("D17").Value = "MAYCA CONSTRUCOES E INSTALACOES ELETROMECANICAS EIRELI"
$ws.Range("E17").Value = "001261"
$ws.Range("F17").Value = "SABAO LIQUIDO LAVA ROUPA BRINORT CONC AZUL 2L"
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = $false
